$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = "56.040.93"
$cell.Style = "Normal"

$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = "2.356.06"
$cell.Style = "Normal"
$ws.Range("E3").Value = "  -1.32%  "

$ws.Range("E4").Value = "  +0.03%  "

$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "505.92"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.26%  "

$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "130.49"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -1.54%  "

$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.997"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E8").Value = "  -2.25%  "

$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "2.371.92"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -0.83%  "

$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "0.0976"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +0.26%  "

$ws.Range("E11").Value = "  -0.65%  "

$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "4.81"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +3.53%  "

$ws.Range("E13").Value = "  -0.71%  "

$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "2.775.52"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -1.22%  "

$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "55.986.79"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("E16").Value = "  -0.46%  "

$ws.Range("E17").Value = "  -0.14%  "

$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = "2.360.78"
$cell.Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "

$ws.Range("E19").Value = "  -2.07%  "

$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "311.89"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +0.96%  "

$ws.Range("E21").Value = "  -0.49%  "

$ws.Range("E23").Value = "  -0.17%  "

$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "65.40"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -0.74%  "

$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("E26").Value = "  -2.71%  "

$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("E28").Value = "  -3.40%  "

$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "171.93"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -2.50%  "

$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "0.0₃0708"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -2.28%  "

$ws.Range("E31").Value = "  -0.32%  "

$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "5.77"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -1.33%  "

$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.996"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.16%  "

$ws.Range("E35").Value = "  -3.77%  "

$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "17.70"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.36%  "

$ws.Range("E37").Value = "  -1.28%  "

$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "0.847"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  +3.04%  "

$ws.Range("E39").Value = "  -3.84%  "

$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "36.18"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -1.45%  "

$ws.Range("E41").Value = "  -3.15%  "

$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "3.36"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -0.67%  "

$ws.Range("E43").Value = "  +0.78%  "

$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "125.60"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -3.99%  "

$ws.Range("E45").Value = "  -1.33%  "

$ws.Range("E46").Value = "  -1.42%  "

$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "242.89"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -2.14%  "

$ws.Range("E48").Value = "  -0.92%  "

$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "16.82"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -0.76%  "

$ws.Range("E50").Value = "  -1.23%  "

$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "16.73"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -2.75%  "
